$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the end-time for the Tuesday (2/26/19) 6:00 PM entry
$ws.Range("A9").Value = "Tuesday (2/26/19) 6:00 PM - 3:00 AM"

# Fill in the rest of that row: Activity, Progress (TODO), Members
$ws.Range("B9").Value = "Fixed debugging issues, ran PRPG program successfully"
$ws.Range("C9").Value = "Test the decoder using Project 1 code part 2"
$ws.Range("D9").Value = "Richard"

# Match styling of the rest of the table rows (center aligned)
$ws.Range("A9:D9").HorizontalAlignment = -4108

# Update the active selection to D9, matching the final cursor position
$ws.Range("D9").Select()
